# Adds the "Enquiry" columns (Z, AA) to Sheet1 row 2 of the Q0016 param sheet.
# Z2 = new long description text ("Premium Paying Terms for Basic Coverage and Riders")
# AA2 = re-uses the existing "Premium Term" text already used in J2.
# Both new cells get a small (7.5pt) wrapped, vertically-centred, thin-bordered style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Values (order matters: Z2 first creates the new shared string,
#     AA2 second re-uses the pre-existing "Premium Term" shared string) ---
$ws.Range("Z2").Value = "Premium Paying Terms for Basic Coverage and Riders"
$ws.Range("AA2").Value = "Premium Term"

# --- Formatting: build the style once on an unused scratch cell, then copy
#     it onto the target cells as a single "paste formats" operation so we
#     do not leave behind a trail of intermediate/unused cell styles. ---
$scratch = $ws.Range("ZZ999")
$scratch.Font.Size = 7.5
$scratch.Borders.LineStyle = 1
$scratch.WrapText = $true
$scratch.VerticalAlignment = -4108  # xlCenter

$target = $ws.Range("Z2:AA2")
$scratch.Copy()
$target.PasteSpecial(-4122)  # xlPasteFormats
$scratch.Clear()

# --- Column sizing for the two new columns ---
$ws.Columns.Item(26).ColumnWidth = 46.5
$ws.Columns.Item(27).ColumnWidth = 50

# --- View state: scroll right a bit and select the newly added cells ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 17
$ws.Range("Z2:AA2").Select()
